# The original "max" column (C) is dropped entirely, shifting the old
# "prediction" (D) and "rejection-f" (E) columns left by one, and the
# "1-g__UBA3663" values in column B become real prediction scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "max" column (C) - this shifts D->C and E->D,
# and also updates the sheet dimension from A1:E3 to A1:D3 automatically.
$ws.Range("C:C").EntireColumn.Delete()

# Update the numeric values in column B (now the "1-g__UBA3663" score column).
$ws.Range("B2").Value = 156.8500202627738
$ws.Range("B3").Value = 206.1248448436627
